# Add materials for session 08 (row 9: week 8)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Column E (Folien) and F (Aufgaben) for session 08 / row 9
# Set F9 first so the shared-string table order matches the authored file
# (exercises/e08.html = index 44, slides link = index 45)
$ws.Range("F9").Value = "exercises/e08.html"
$ws.Range("E9").Value = "slides/slides.html#/sitzung-08-codebuch-goldstandard"

# Update the active cell selection to match the authored state
$ws.Range("D9").Select()
